$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The workbook had two "_xlnm._FilterDatabase" defined names (one hidden,
#    one visible/stale). Rename the visible one so it no longer collides
#    with the hidden autofilter name.
foreach ($n in @($wb.Names)) {
    if ($n.Visible -and $n.Name -eq "Sheet1!_FilterDatabase") {
        $n.Name = "_xlnm._FilterDatabase_0"
    }
}

# 2. The "Needs Coverage" doctype values used the old tag names; switch the
#    whole column over to the new built-in default doctype names.
$ws.Range("E2:E76").Replace("impl;test", "sourcecode;testcode") | Out-Null

# 3. Move the active selection.
$ws.Range("A7").Select() | Out-Null

# 4. Minor column C width touch-up recorded by the spreadsheet app.
$ws.Columns.Item(3).ColumnWidth = 10.3
